$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) cells: force text format so values like "218.37" or "1.01"
# are not reinterpreted as numbers, matching the inline-string text cells in the source.
$dCells = @("D2", "D3", "D5", "D9", "D10", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D31", "D32", "D35", "D37", "D38", "D42", "D43", "D45", "D46", "D48", "D50", "D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.796.84"
$ws.Range("D3").Value = "1.639.96"
$ws.Range("D5").Value = "218.37"
$ws.Range("D9").Value = "0.0621"
$ws.Range("D10").Value = "19.28"
$ws.Range("D12").Value = "1.868.28"
$ws.Range("D13").Value = "1.634.38"
$ws.Range("D14").Value = "4.14"
$ws.Range("D15").Value = "0.525"
$ws.Range("D16").Value = "64.99"
$ws.Range("D17").Value = "26.788.06"
$ws.Range("D18").Value = "0.0₃0731"
$ws.Range("D19").Value = "216.22"
$ws.Range("D23").Value = "2.35"
$ws.Range("D24").Value = "9.15"
$ws.Range("D25").Value = "146.98"
$ws.Range("D26").Value = "1.01"
$ws.Range("D27").Value = "0.118"
$ws.Range("D28").Value = "7.09"
$ws.Range("D29").Value = "15.73"
$ws.Range("D31").Value = "1.19"
$ws.Range("D32").Value = "3.38"
$ws.Range("D35").Value = "1.262.49"
$ws.Range("D37").Value = "0.0174"
$ws.Range("D38").Value = "0.531"
$ws.Range("D42").Value = "5.35"
$ws.Range("D43").Value = "1.780.25"
$ws.Range("D45").Value = "92.20"
$ws.Range("D46").Value = "60.81"
$ws.Range("D48").Value = "0.0₆0103"
$ws.Range("D50").Value = "7.53"
$ws.Range("D51").Value = "0.0962"

foreach ($addr in $dCells) {
    $ws.Range($addr).ClearFormats()
}

# Column E (Volume 1h) cells: plain text assignment (percentage strings with padding
# are already safely treated as text by Excel).
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("E7").Value = "  -0.54%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -0.84%  "
$ws.Range("E10").Value = "  +0.48%  "
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("E13").Value = "  -0.89%  "
$ws.Range("E14").Value = "  -0.70%  "
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("E16").Value = "  +0.60%  "
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("E20").Value = "  -0.38%  "
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("E22").Value = "  +4.66%  "
$ws.Range("E23").Value = "  -1.97%  "
$ws.Range("E24").Value = "  -2.27%  "
$ws.Range("E25").Value = "  +1.32%  "
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("E30").Value = "  -1.80%  "
$ws.Range("E31").Value = "  +0.96%  "
$ws.Range("E32").Value = "  +1.80%  "
$ws.Range("E33").Value = "  -0.78%  "
$ws.Range("E34").Value = "  +0.65%  "
$ws.Range("E35").Value = "  -2.05%  "
$ws.Range("E36").Value = "  +0.35%  "
$ws.Range("E37").Value = "  -0.76%  "
$ws.Range("E38").Value = "  -1.38%  "
$ws.Range("E39").Value = "  -1.25%  "
$ws.Range("E40").Value = "  -0.33%  "
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("E42").Value = "  -0.36%  "
$ws.Range("E44").Value = "  -4.21%  "
$ws.Range("E45").Value = "  +0.83%  "
$ws.Range("E46").Value = "  +0.71%  "
$ws.Range("E47").Value = "  -1.13%  "
$ws.Range("E48").Value = "  -2.84%  "
$ws.Range("E49").Value = "  -0.69%  "
$ws.Range("E50").Value = "  -1.71%  "
$ws.Range("E51").Value = "  -1.84%  "
